# Applies the "Updated symbol list" price/coin refresh described by the diff.
# Column D (Price) cells are stored as literal TEXT in this sheet (not numbers),
# so numeric-looking updates are written with a leading apostrophe to force text
# storage (preserving exact digits/trailing zeros) and then the cell Style is reset
# to "Normal" so the quote-prefix formatting flag does not leak into the output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'249.12"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'22.68"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Value = "'5.239"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'0.05686"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Value = "'6.343"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'0.8049"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.9173"
$ws.Range('D9').Style = 'Normal'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = "'0.1406"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = "'0.07421"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = "'0.03108"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03030"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09377"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').Value = "'3.891"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = "'0.001587"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = "'0.04792"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'UpBots'
$ws.Range('C18').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D18').Value = "'0.01828"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '17UpBotsUBXTBestin24h'
$ws.Range('B19').Value = 'One'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D19').Value = "'0.0005852"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '18OneONE'
$ws.Range('D21').Value = "'0.004993"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'0.0009992"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'0.0001500"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'3.696"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Value = "'2.195"
$ws.Range('D25').Style = 'Normal'
$ws.Range('D27').Value = "'0.1307"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D40').Value = "'0.04001"
$ws.Range('D40').Style = 'Normal'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = "'0.1071"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = "'0.002740"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = "'0.006750"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('D44').Value = "'0.007960"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'0.00005661"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range('D48').Value = "'0.2083"
$ws.Range('D48').Style = 'Normal'
